$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Header strings: bump the report Volume/Number and the reporting week
#    dates. These live as multi-run rich text inside a single shared-string
#    cell each, so we surgically replace just the affected run's characters
#    to keep the rest of the cell text (and its other runs) intact.
# ---------------------------------------------------------------------------

# "Volume 31   Number  43" -> "...44"   (the "43" run starts at char 21, len 2)
$ws.Range("A8").Characters(21, 2).Text = "44"

# "Report Covering the Week  10/21/2024  Through  10/27/2024"
#   -> "...10/28/2024  Through  11/3/2024"
$ws.Range("C9").Characters(27, 10).Text = "10/28/2024"
$ws.Range("C9").Characters(48, 10).Text = "11/3/2024"

# ---------------------------------------------------------------------------
# 2. Plain numeric value updates (weekly crime-stat counts/percentages).
#    These cells keep their existing style/number-format, only the value
#    itself changes, so a simple .Value assignment is enough.
# ---------------------------------------------------------------------------
$pairs = @{
    "F16" = 3
    "G16" = 2
    "H16" = 50
    "I16" = 50
    "K16" = 56.25
    "L16" = 8.695652173913
    "M16" = -13.793103448275
    "N16" = -78.902953586497
    "D17" = 2
    "E17" = 0
    "F17" = 8
    "G17" = 16
    "H17" = -50
    "I17" = 135
    "J17" = 111
    "K17" = 21.621621621621
    "L17" = 48.351648351648
    "M17" = 84.931506849315
    "N17" = -29.319371727748
    "G18" = 4
    "H18" = 50
    "N18" = -86.849315068493
    "D19" = 3
    "E19" = -66.666666666666
    "F19" = 6
    "H19" = -64.705882352941
    "I19" = 133
    "J19" = 162
    "K19" = -17.901234567901
    "L19" = -5
    "M19" = 31.683168316831
    "N19" = -23.121387283237
    "G20" = 6
    "H20" = -83.333333333333
    "I20" = 48
    "K20" = 33.333333333333
    "L20" = -12.727272727272
    "M20" = 92
    "N20" = -83.216783216783
    "C21" = 5
    "D21" = 5
    "E21" = 0
    "F21" = 24
    "G21" = 46
    "H21" = -47.826086956521
    "I21" = 421
    "J21" = 381
    "K21" = 10.498687664042
    "L21" = 7.124681933842
    "M21" = 35.806451612903
    "N21" = -66.824271079590
    "L22" = 44.444444444444
    "C24" = 3
    "D24" = 6
    "E24" = -50
    "F24" = 39
    "G24" = 35
    "H24" = 11.428571428571
    "I24" = 429
    "J24" = 411
    "K24" = 4.379562043795
    "L24" = -5.088495575221
    "M24" = 65.637065637065
    "C25" = 4
    "D25" = 2
    "E25" = 100
    "G25" = 10
    "H25" = 200
    "I25" = 223
    "J25" = 213
    "K25" = 4.694835680751
    "L25" = 12.060301507537
    "C26" = 5
    "D26" = 3
    "E26" = 66.666666666666
    "F26" = 11
    "H26" = -31.25
    "I26" = 226
    "J26" = 190
    "K26" = 18.947368421052
    "L26" = 22.826086956521
    "M26" = -22.602739726027
    "J28" = 19
    "K28" = -10.526315789473
    "L28" = -15
}

foreach ($ref in $pairs.Keys) {
    $ws.Range($ref).Value = $pairs[$ref]
}

# ---------------------------------------------------------------------------
# 3. Cells that flip from a numeric count to the "no data" text placeholders
#    ("0" / "***.*", shared strings already present in the workbook). Excel
#    auto-detects digit-only strings as numbers unless the cell is flagged
#    as Text first, so force the Text number format for the assignment and
#    then drop back to General (matches the target cell's rendering).
# ---------------------------------------------------------------------------
function Set-TextPlaceholder($ref, $text) {
    $ws.Range($ref).NumberFormat = "@"
    $ws.Range($ref).Value = $text
    $ws.Range($ref).NumberFormat = "General"
}

Set-TextPlaceholder "D15" "0"
Set-TextPlaceholder "E15" "***.*"
Set-TextPlaceholder "D16" "0"
Set-TextPlaceholder "E16" "***.*"
Set-TextPlaceholder "C18" "0"
Set-TextPlaceholder "D18" "0"
Set-TextPlaceholder "E18" "***.*"
Set-TextPlaceholder "D20" "0"
Set-TextPlaceholder "E20" "***.*"
Set-TextPlaceholder "D27" "0"
Set-TextPlaceholder "E27" "***.*"
Set-TextPlaceholder "C28" "0"

# ---------------------------------------------------------------------------
# 4. Cells that flip the other way: from the "no data" text placeholder back
#    to a real numeric count/percentage. Re-apply the destination cell's
#    original numeric format explicitly so the underlying style matches the
#    numeric sibling cells (format 167 = "#,##0" / format 166 = count-style
#    decimal used for the percentage columns) rather than creating a new
#    ad-hoc style.
# ---------------------------------------------------------------------------
$ws.Range("C17").NumberFormat = "#,##0"
$ws.Range("C17").Value = 2

$ws.Range("C20").NumberFormat = "#,##0"
$ws.Range("C20").Value = 1

$ws.Range("F20").NumberFormat = "#,##0"
$ws.Range("F20").Value = 1

$ws.Range("D28").NumberFormat = "#,##0"
$ws.Range("D28").Value = 1

$ws.Range("E28").NumberFormat = '#,##0.0;"-"#,##0.0'
$ws.Range("E28").Value = -100

$ws.Range("G28").NumberFormat = "#,##0"
$ws.Range("G28").Value = 1

$ws.Range("H28").NumberFormat = '#,##0.0;"-"#,##0.0'
$ws.Range("H28").Value = 0
